# Add data for 2024-06-04
# Updates the 2024 (column K) violent-crime counts across the citywide
# totals sheet, the by-neighborhood summary sheet, and every individual
# neighborhood sheet, reflecting one additional day of incident data.
# A handful of prior-year cells (2017/2021/2022/2023) also receive small
# corrections that came in alongside this day's data.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Cells.Item(2, 11).Value = 3193  # K2: 3162 -> 3193
$ws.Cells.Item(3, 11).Value = 3168  # K3: 3140 -> 3168
$ws.Cells.Item(4, 4).Value = 1974  # D4: 1973 -> 1974
$ws.Cells.Item(4, 8).Value = 1730  # H4: 1729 -> 1730
$ws.Cells.Item(4, 9).Value = 1794  # I4: 1795 -> 1794
$ws.Cells.Item(4, 10).Value = 1817  # J4: 1818 -> 1817
$ws.Cells.Item(4, 11).Value = 652  # K4: 645 -> 652
$ws.Cells.Item(5, 11).Value = 210  # K5: 207 -> 210
$ws.Cells.Item(6, 11).Value = 3753  # K6: 3718 -> 3753
$ws.Cells.Item(7, 4).Value = 28164  # D7: 28163 -> 28164
$ws.Cells.Item(7, 8).Value = 26043  # H7: 26042 -> 26043
$ws.Cells.Item(7, 9).Value = 26248  # I7: 26249 -> 26248
$ws.Cells.Item(7, 10).Value = 29288  # J7: 29289 -> 29288
$ws.Cells.Item(7, 11).Value = 10976  # K7: 10872 -> 10976

$ws = $wb.Worksheets.Item("Logan Square")
$ws.Cells.Item(4, 11).Value = 6  # K4: 5 -> 6
$ws.Cells.Item(6, 11).Value = 77  # K6: 76 -> 77
$ws.Cells.Item(7, 11).Value = 149  # K7: 147 -> 149

$ws = $wb.Worksheets.Item("Austin")
$ws.Cells.Item(2, 11).Value = 211  # K2: 207 -> 211
$ws.Cells.Item(3, 11).Value = 223  # K3: 221 -> 223
$ws.Cells.Item(4, 11).Value = 40  # K4: 39 -> 40
$ws.Cells.Item(5, 11).Value = 20  # K5: 19 -> 20
$ws.Cells.Item(6, 11).Value = 238  # K6: 235 -> 238
$ws.Cells.Item(7, 11).Value = 732  # K7: 721 -> 732

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Cells.Item(6, 11).Value = 55  # K6: 54 -> 55
$ws.Cells.Item(7, 11).Value = 235  # K7: 234 -> 235

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Cells.Item(2, 11).Value = 120  # K2: 119 -> 120
$ws.Cells.Item(3, 11).Value = 155  # K3: 154 -> 155
$ws.Cells.Item(7, 11).Value = 429  # K7: 427 -> 429

$ws = $wb.Worksheets.Item("West Pullman")
$ws.Cells.Item(2, 11).Value = 58  # K2: 57 -> 58
$ws.Cells.Item(7, 11).Value = 179  # K7: 178 -> 179

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Cells.Item(2, 11).Value = 104  # K2: 101 -> 104
$ws.Cells.Item(3, 11).Value = 130  # K3: 128 -> 130
$ws.Cells.Item(6, 11).Value = 116  # K6: 115 -> 116
$ws.Cells.Item(7, 11).Value = 381  # K7: 375 -> 381

$ws = $wb.Worksheets.Item("New City")
$ws.Cells.Item(3, 11).Value = 66  # K3: 64 -> 66
$ws.Cells.Item(6, 11).Value = 101  # K6: 100 -> 101
$ws.Cells.Item(7, 11).Value = 258  # K7: 255 -> 258

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Cells.Item(2, 11).Value = 86  # K2: 85 -> 86
$ws.Cells.Item(6, 8).Value = 171  # H6: 170 -> 171
$ws.Cells.Item(6, 11).Value = 85  # K6: 84 -> 85
$ws.Cells.Item(8, 11).Value = 732  # K8: 721 -> 732
$ws.Cells.Item(11, 11).Value = 235  # K11: 232 -> 235
$ws.Cells.Item(13, 11).Value = 13  # K13: 12 -> 13
$ws.Cells.Item(15, 11).Value = 113  # K15: 112 -> 113
$ws.Cells.Item(18, 11).Value = 77  # K18: 76 -> 77
$ws.Cells.Item(19, 11).Value = 332  # K19: 329 -> 332
$ws.Cells.Item(20, 11).Value = 252  # K20: 248 -> 252
$ws.Cells.Item(23, 11).Value = 102  # K23: 101 -> 102
$ws.Cells.Item(27, 11).Value = 112  # K27: 111 -> 112
$ws.Cells.Item(29, 11).Value = 574  # K29: 566 -> 574
$ws.Cells.Item(31, 11).Value = 120  # K31: 119 -> 120
$ws.Cells.Item(33, 11).Value = 429  # K33: 427 -> 429
$ws.Cells.Item(36, 11).Value = 130  # K36: 126 -> 130
$ws.Cells.Item(37, 11).Value = 381  # K37: 375 -> 381
$ws.Cells.Item(41, 11).Value = 97  # K41: 95 -> 97
$ws.Cells.Item(42, 9).Value = 1009  # I42: 1010 -> 1009
$ws.Cells.Item(42, 11).Value = 381  # K42: 375 -> 381
$ws.Cells.Item(43, 11).Value = 97  # K43: 96 -> 97
$ws.Cells.Item(48, 11).Value = 134  # K48: 133 -> 134
$ws.Cells.Item(51, 11).Value = 127  # K51: 122 -> 127
$ws.Cells.Item(52, 11).Value = 299  # K52: 298 -> 299
$ws.Cells.Item(53, 11).Value = 149  # K53: 147 -> 149
$ws.Cells.Item(54, 11).Value = 212  # K54: 211 -> 212
$ws.Cells.Item(55, 11).Value = 118  # K55: 115 -> 118
$ws.Cells.Item(63, 4).Value = 353  # D63: 352 -> 353
$ws.Cells.Item(63, 10).Value = 101  # J63: 102 -> 101
$ws.Cells.Item(63, 11).Value = 40  # K63: 38 -> 40
$ws.Cells.Item(65, 11).Value = 258  # K65: 255 -> 258
$ws.Cells.Item(67, 11).Value = 433  # K67: 430 -> 433
$ws.Cells.Item(75, 11).Value = 39  # K75: 38 -> 39
$ws.Cells.Item(76, 11).Value = 166  # K76: 165 -> 166
$ws.Cells.Item(77, 11).Value = 78  # K77: 77 -> 78
$ws.Cells.Item(79, 11).Value = 283  # K79: 279 -> 283
$ws.Cells.Item(80, 11).Value = 36  # K80: 35 -> 36
$ws.Cells.Item(83, 11).Value = 235  # K83: 234 -> 235
$ws.Cells.Item(85, 11).Value = 521  # K85: 513 -> 521
$ws.Cells.Item(88, 11).Value = 131  # K88: 128 -> 131
$ws.Cells.Item(89, 11).Value = 144  # K89: 142 -> 144
$ws.Cells.Item(90, 11).Value = 97  # K90: 95 -> 97
$ws.Cells.Item(91, 11).Value = 115  # K91: 114 -> 115
$ws.Cells.Item(95, 11).Value = 179  # K95: 178 -> 179
$ws.Cells.Item(98, 11).Value = 60  # K98: 61 -> 60
$ws.Cells.Item(101, 4).Value = 28164  # D101: 28163 -> 28164
$ws.Cells.Item(101, 8).Value = 26043  # H101: 26042 -> 26043
$ws.Cells.Item(101, 9).Value = 26248  # I101: 26249 -> 26248
$ws.Cells.Item(101, 10).Value = 29288  # J101: 29289 -> 29288
$ws.Cells.Item(101, 11).Value = 10976  # K101: 10872 -> 10976

$ws = $wb.Worksheets.Item("Gage Park")
$ws.Cells.Item(3, 11).Value = 27  # K3: 26 -> 27
$ws.Cells.Item(7, 11).Value = 120  # K7: 119 -> 120

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Cells.Item(2, 11).Value = 135  # K2: 134 -> 135
$ws.Cells.Item(6, 11).Value = 124  # K6: 122 -> 124
$ws.Cells.Item(7, 11).Value = 433  # K7: 430 -> 433

$ws = $wb.Worksheets.Item("Loop")
$ws.Cells.Item(6, 11).Value = 97  # K6: 96 -> 97
$ws.Cells.Item(7, 11).Value = 212  # K7: 211 -> 212

$ws = $wb.Worksheets.Item("Englewood")
$ws.Cells.Item(2, 11).Value = 157  # K2: 153 -> 157
$ws.Cells.Item(3, 11).Value = 200  # K3: 198 -> 200
$ws.Cells.Item(4, 11).Value = 32  # K4: 31 -> 32
$ws.Cells.Item(6, 11).Value = 173  # K6: 172 -> 173
$ws.Cells.Item(7, 11).Value = 574  # K7: 566 -> 574

$ws = $wb.Worksheets.Item("Lake View")
$ws.Cells.Item(4, 11).Value = 18  # K4: 17 -> 18
$ws.Cells.Item(7, 11).Value = 134  # K7: 133 -> 134

$ws = $wb.Worksheets.Item("Chatham")
$ws.Cells.Item(2, 11).Value = 112  # K2: 111 -> 112
$ws.Cells.Item(3, 11).Value = 88  # K3: 86 -> 88
$ws.Cells.Item(7, 11).Value = 332  # K7: 329 -> 332

$ws = $wb.Worksheets.Item("River North")
$ws.Cells.Item(6, 11).Value = 98  # K6: 97 -> 98
$ws.Cells.Item(7, 11).Value = 166  # K7: 165 -> 166

$ws = $wb.Worksheets.Item("Ashburn")
$ws.Cells.Item(2, 11).Value = 32  # K2: 31 -> 32
$ws.Cells.Item(4, 8).Value = 9  # H4: 8 -> 9
$ws.Cells.Item(7, 8).Value = 171  # H7: 170 -> 171
$ws.Cells.Item(7, 11).Value = 85  # K7: 84 -> 85

$ws = $wb.Worksheets.Item("Hermosa")
$ws.Cells.Item(4, 11).Value = 6  # K4: 5 -> 6
$ws.Cells.Item(6, 11).Value = 41  # K6: 40 -> 41
$ws.Cells.Item(7, 11).Value = 97  # K7: 95 -> 97

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Cells.Item(3, 11).Value = 119  # K3: 118 -> 119
$ws.Cells.Item(4, 9).Value = 56  # I4: 57 -> 56
$ws.Cells.Item(6, 11).Value = 147  # K6: 142 -> 147
$ws.Cells.Item(7, 9).Value = 1009  # I7: 1010 -> 1009
$ws.Cells.Item(7, 11).Value = 381  # K7: 375 -> 381

$ws = $wb.Worksheets.Item("Boystown")
$ws.Cells.Item(3, 11).Value = 4  # K3: 3 -> 4
$ws.Cells.Item(6, 11).Value = 13  # K6: 12 -> 13

$ws = $wb.Worksheets.Item("Lower West Side")
$ws.Cells.Item(3, 11).Value = 29  # K3: 27 -> 29
$ws.Cells.Item(6, 11).Value = 43  # K6: 42 -> 43
$ws.Cells.Item(7, 11).Value = 118  # K7: 115 -> 118

$ws = $wb.Worksheets.Item("Douglas")
$ws.Cells.Item(6, 11).Value = 27  # K6: 26 -> 27
$ws.Cells.Item(7, 11).Value = 102  # K7: 101 -> 102

$ws = $wb.Worksheets.Item("Washington Park")
$ws.Cells.Item(3, 11).Value = 55  # K3: 54 -> 55
$ws.Cells.Item(7, 11).Value = 115  # K7: 114 -> 115

$ws = $wb.Worksheets.Item("Roseland")
$ws.Cells.Item(2, 11).Value = 99  # K2: 97 -> 99
$ws.Cells.Item(6, 11).Value = 62  # K6: 60 -> 62
$ws.Cells.Item(7, 11).Value = 283  # K7: 279 -> 283

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Cells.Item(2, 11).Value = 86  # K2: 84 -> 86
$ws.Cells.Item(3, 11).Value = 70  # K3: 69 -> 70
$ws.Cells.Item(6, 11).Value = 84  # K6: 83 -> 84
$ws.Cells.Item(7, 11).Value = 252  # K7: 248 -> 252

$ws = $wb.Worksheets.Item("Calumet Heights")
$ws.Cells.Item(3, 11).Value = 25  # K3: 24 -> 25
$ws.Cells.Item(7, 11).Value = 77  # K7: 76 -> 77

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Cells.Item(2, 11).Value = 51  # K2: 50 -> 51
$ws.Cells.Item(3, 11).Value = 39  # K3: 38 -> 39
$ws.Cells.Item(4, 11).Value = 13  # K4: 12 -> 13
$ws.Cells.Item(6, 11).Value = 27  # K6: 26 -> 27
$ws.Cells.Item(7, 11).Value = 130  # K7: 126 -> 130

$ws = $wb.Worksheets.Item("Brighton Park")
$ws.Cells.Item(2, 11).Value = 40  # K2: 39 -> 40
$ws.Cells.Item(7, 11).Value = 113  # K7: 112 -> 113

$ws = $wb.Worksheets.Item("Wicker Park")
$ws.Cells.Item(6, 11).Value = 44  # K6: 45 -> 44
$ws.Cells.Item(7, 11).Value = 60  # K7: 61 -> 60

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Cells.Item(2, 11).Value = 69  # K2: 68 -> 69
$ws.Cells.Item(6, 11).Value = 90  # K6: 88 -> 90
$ws.Cells.Item(7, 11).Value = 235  # K7: 232 -> 235

$ws = $wb.Worksheets.Item("Albany Park")
$ws.Cells.Item(2, 11).Value = 23  # K2: 22 -> 23
$ws.Cells.Item(7, 11).Value = 86  # K7: 85 -> 86

$ws = $wb.Worksheets.Item("United Center")
$ws.Cells.Item(2, 11).Value = 29  # K2: 28 -> 29
$ws.Cells.Item(3, 11).Value = 33  # K3: 31 -> 33
$ws.Cells.Item(7, 11).Value = 131  # K7: 128 -> 131

$ws = $wb.Worksheets.Item("Uptown")
$ws.Cells.Item(4, 11).Value = 21  # K4: 20 -> 21
$ws.Cells.Item(6, 11).Value = 42  # K6: 41 -> 42
$ws.Cells.Item(7, 11).Value = 144  # K7: 142 -> 144

$ws = $wb.Worksheets.Item("Edgewater")
$ws.Cells.Item(2, 11).Value = 30  # K2: 29 -> 30
$ws.Cells.Item(7, 11).Value = 112  # K7: 111 -> 112

$ws = $wb.Worksheets.Item("Pullman")
$ws.Cells.Item(3, 11).Value = 11  # K3: 10 -> 11
$ws.Cells.Item(7, 11).Value = 39  # K7: 38 -> 39

$ws = $wb.Worksheets.Item("Washington Heights")
$ws.Cells.Item(3, 11).Value = 31  # K3: 30 -> 31
$ws.Cells.Item(5, 11).Value = 3  # K5: 2 -> 3
$ws.Cells.Item(7, 11).Value = 97  # K7: 95 -> 97

$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Cells.Item(2, 11).Value = 34  # K2: 32 -> 34
$ws.Cells.Item(3, 11).Value = 31  # K3: 30 -> 31
$ws.Cells.Item(6, 11).Value = 48  # K6: 46 -> 48
$ws.Cells.Item(7, 11).Value = 127  # K7: 122 -> 127

$ws = $wb.Worksheets.Item("Hyde Park")
$ws.Cells.Item(3, 11).Value = 29  # K3: 28 -> 29
$ws.Cells.Item(7, 11).Value = 97  # K7: 96 -> 97

$ws = $wb.Worksheets.Item("South Shore")
$ws.Cells.Item(3, 11).Value = 179  # K3: 177 -> 179
$ws.Cells.Item(4, 11).Value = 27  # K4: 26 -> 27
$ws.Cells.Item(5, 11).Value = 12  # K5: 11 -> 12
$ws.Cells.Item(6, 11).Value = 118  # K6: 114 -> 118
$ws.Cells.Item(7, 11).Value = 521  # K7: 513 -> 521

$ws = $wb.Worksheets.Item("Riverdale")
$ws.Cells.Item(6, 11).Value = 8  # K6: 7 -> 8
$ws.Cells.Item(7, 11).Value = 78  # K7: 77 -> 78

$ws = $wb.Worksheets.Item("Rush & Division")
$ws.Cells.Item(2, 11).Value = 9  # K2: 8 -> 9
$ws.Cells.Item(7, 11).Value = 36  # K7: 35 -> 36

$ws = $wb.Worksheets.Item("Little Village")
$ws.Cells.Item(6, 11).Value = 120  # K6: 119 -> 120
$ws.Cells.Item(7, 11).Value = 299  # K7: 298 -> 299
